$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New week-37 header cell (AN1): must stay text "37" like the other week headers ---
$ws.Range("AN1").Value = "'37"
$ws.Range("AM1").Copy() | Out-Null
$ws.Range("AN1").PasteSpecial(-4122) | Out-Null

# --- Data cell updates (new AN column values + corrected week figures) ---
$ws.Range("AN2").Value = 0
$ws.Range("AN3").Value = 0
$ws.Range("AN5").Value = 0
$ws.Range("D6").Value = 2
$ws.Range("I6").Value = 2
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 2
$ws.Range("U6").Value = 3
$ws.Range("W6").Value = 2
$ws.Range("Y6").Value = 1
$ws.Range("Z6").Value = 4
$ws.Range("AA6").Value = 3
$ws.Range("AB6").Value = 4
$ws.Range("AC6").Value = 1
$ws.Range("AD6").Value = 2
$ws.Range("AE6").Value = 4
$ws.Range("AF6").Value = 4
$ws.Range("AG6").Value = 3
$ws.Range("AH6").Value = 3
$ws.Range("AN6").Value = 1
$ws.Range("AN7").Value = 0
$ws.Range("AN8").Value = 0
$ws.Range("AN9").Value = 0
$ws.Range("AN10").Value = 0
$ws.Range("AN12").Value = 0
$ws.Range("AN14").Value = 0
$ws.Range("AN15").Value = 0
$ws.Range("AN16").Value = 0
$ws.Range("AN17").Value = 0
$ws.Range("AN21").Value = 0
$ws.Range("AN23").Value = 0
$ws.Range("AN25").Value = 0
$ws.Range("AN26").Value = 0
$ws.Range("T28").Value = 1
$ws.Range("U28").Value = 1
$ws.Range("W28").Value = 0
$ws.Range("X28").Value = 0
$ws.Range("Z28").Value = 0
$ws.Range("AE28").Value = 1
$ws.Range("AF28").Value = 0
$ws.Range("AG28").Value = 2
$ws.Range("AH28").Value = 1
$ws.Range("AI28").Value = 0
$ws.Range("AJ28").Value = 0
$ws.Range("AK28").Value = 0
$ws.Range("AN28").Value = 0
$ws.Range("AN29").Value = 1
$ws.Range("AN30").Value = 4
$ws.Range("AM31").Value = 0
$ws.Range("AN31").Value = 0
$ws.Range("M35").Value = 2
$ws.Range("AM35").Value = 1
$ws.Range("AN35").Value = 3
$ws.Range("AN36").Value = 0
$ws.Range("AN37").Value = 0
$ws.Range("AN38").Value = 0
$ws.Range("AN41").Value = 0
$ws.Range("AN42").Value = 0
$ws.Range("AM43").Value = 0
$ws.Range("AN43").Value = 0
$ws.Range("AN44").Value = 0
$ws.Range("AN45").Value = 0
$ws.Range("AN46").Value = 0
$ws.Range("AN47").Value = 0
$ws.Range("AN48").Value = 0
$ws.Range("AN49").Value = 0
$ws.Range("AN50").Value = 0
$ws.Range("AN51").Value = 0
$ws.Range("AM52").Value = 0
$ws.Range("AN53").Value = 0
$ws.Range("AN54").Value = 0
$ws.Range("AN55").Value = 0
$ws.Range("AN56").Value = 0
$ws.Range("AN57").Value = 0
$ws.Range("AN58").Value = 0
